# "support forward iterator for bubble sort & odd even sort"
#
# The "Summary" sheet lists sorting algorithms with an "Iterator" column
# (F). Bubble Sort (row 2) and Odd Even Sort (row 5) previously required a
# bidirectional iterator; they now only need a forward iterator.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Summary")

$ws.Range("F2").Value = "forward iterator"   # Bubble Sort
$ws.Range("F5").Value = "forward iterator"   # Odd Even Sort

# Leave the selection where the author ended up after making the edit.
$ws.Range("F7").Select()
